$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.04983966666666667
$ws.Range("H2").Value = 0.149519
$ws.Range("I2").Value = 0.1823731600337622
$ws.Range("J2").Value = 0.1823731600337622
$ws.Range("M2").Value = 0.4102596666666667
$ws.Range("N2").Value = 1.230779
$ws.Range("O2").Value = 0.003499619873322347
$ws.Range("P2").Value = 0.003499619873322347
$ws.Range("Q2").Value = 0.02044720503344445
$ws.Range("R2").Value = 0.184024845301
$ws.Range("S2").Value = 0.0006382367352147511
$ws.Range("T2").Value = 0.000638236735214751
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.04983966666666667
$ws.Range("H3").Value = 0.149519
$ws.Range("I3").Value = 0.1823731600337622
$ws.Range("J3").Value = 0.1823731600337622
$ws.Range("O3").Value = 0.8692174743460166
$ws.Range("P3").Value = 0.8692174743460165
$ws.Range("Q3").Value = 5.078570976262334
$ws.Range("R3").Value = 45.70713878636101
$ws.Range("S3").Value = 0.1585219375530487
$ws.Range("T3").Value = 0.1585219375530486
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.04983966666666667
$ws.Range("H4").Value = 0.149519
$ws.Range("I4").Value = 0.1823731600337622
$ws.Range("J4").Value = 0.1823731600337622
$ws.Range("N4").Value = 44.764041
$ws.Range("O4").Value = 0.1272829057806611
$ws.Range("P4").Value = 0.1272829057806611
$ws.Range("Q4").Value = 0.7436749606976667
$ws.Range("R4").Value = 6.693074646279
$ws.Range("S4").Value = 0.02321298574549879
$ws.Range("T4").Value = 0.02321298574549879
$ws.Range("I5").Value = 0.4031422744592926
$ws.Range("J5").Value = 0.4031422744592926
$ws.Range("M5").Value = 0.4102596666666667
$ws.Range("N5").Value = 1.230779
$ws.Range("O5").Value = 0.003499619873322347
$ws.Range("P5").Value = 0.003499619873322347
$ws.Range("Q5").Value = 0.04519926474922221
$ws.Range("R5").Value = 0.406793382743
$ws.Range("S5").Value = 0.001410844715474113
$ws.Range("T5").Value = 0.001410844715474112
$ws.Range("I6").Value = 0.4031422744592926
$ws.Range("J6").Value = 0.4031422744592926
$ws.Range("O6").Value = 0.8692174743460166
$ws.Range("P6").Value = 0.8692174743460165
$ws.Range("S6").Value = 0.350418309607615
$ws.Range("T6").Value = 0.3504183096076149
$ws.Range("I7").Value = 0.4031422744592926
$ws.Range("J7").Value = 0.4031422744592926
$ws.Range("N7").Value = 44.764041
$ws.Range("O7").Value = 0.1272829057806611
$ws.Range("P7").Value = 0.1272829057806611
$ws.Range("S7").Value = 0.05131312013620358
$ws.Range("T7").Value = 0.05131312013620357
$ws.Range("I8").Value = 0.4144845655069452
$ws.Range("J8").Value = 0.4144845655069451
$ws.Range("M8").Value = 0.4102596666666667
$ws.Range("N8").Value = 1.230779
$ws.Range("O8").Value = 0.003499619873322347
$ws.Range("P8").Value = 0.003499619873322347
$ws.Range("Q8").Value = 0.04647093296266667
$ws.Range("R8").Value = 0.418238396664
$ws.Range("S8").Value = 0.001450538422633484
$ws.Range("T8").Value = 0.001450538422633483
$ws.Range("I9").Value = 0.4144845655069452
$ws.Range("J9").Value = 0.4144845655069451
$ws.Range("O9").Value = 0.8692174743460166
$ws.Range("P9").Value = 0.8692174743460165
$ws.Range("S9").Value = 0.3602772271853529
$ws.Range("T9").Value = 0.3602772271853529
$ws.Range("I10").Value = 0.4144845655069452
$ws.Range("J10").Value = 0.4144845655069451
$ws.Range("N10").Value = 44.764041
$ws.Range("O10").Value = 0.1272829057806611
$ws.Range("P10").Value = 0.1272829057806611
$ws.Range("S10").Value = 0.05275679989895878
$ws.Range("T10").Value = 0.05275679989895877
